# Update "想去人数" (interested-count) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 4;  Value = 8319 },
    @{ Row = 5;  Value = 6066 },
    @{ Row = 6;  Value = 522 },
    @{ Row = 7;  Value = 106 },
    @{ Row = 10; Value = 315 }
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}
$ws1.Cells.Item(11, 6).Value = 1012

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Value
}
$ws4.Cells.Item(15, 6).Value = 1012
